# Generate Report for handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets to reflect the newly
# generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 10:56:07"
$wsZhCn.Range("G3").Value = "2016-01-18 10:56:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 10:56:16"
$wsDeDe.Range("G3").Value = "2016-01-18 10:57:10"
